$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.279124021530151
$ws.Range("B1").Value = 2.931780815124512
$ws.Range("C1").Value = 5.378939151763916
$ws.Range("D1").Value = 1.856570720672607
$ws.Range("E1").Value = 1.021705985069275
